$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# --- Sheet1: add new "END_ROW true" / "END_ROW" marker column (J) ---
$ws1.Range("J2").Value = "#! END_ROW true"
$ws1.Range("J3").Value = "#! END_ROW"
$ws1.Range("G18").Select()

# --- Sheet2: add new "END_ROW true" / "END_ROW" marker column (J) ---
$ws2.Range("J2").Value = "#! END_ROW true"
$ws2.Range("J3").Value = "#! END_ROW true"
$ws2.Range("J4").Value = "#! END_ROW"

$ws2.Activate()
$ws2.Range("J4").Select()
